# Consolidate the two header rows (old merged row 1 + detail row 2) into a
# single header row, combining the stakeholder context directly into the
# column labels (so the sheet reads well without merged cells on mobile).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Deleting row 1 (the blank/merged "Stakeholder S (1)/(2)" banner row) shifts
# every row below it up by one: the old detail-header row 2 becomes row 1,
# the data rows 3-17 become rows 2-16, and the summary rows 18-19 become
# rows 17-18. Formulas referencing the shifted rows are adjusted automatically.
$ws.Rows.Item(1).Delete()

# Fold the stakeholder labels (that used to live in the merged banner row)
# into the value/urgency column headers that are now on row 1.
$ws.Range("D1").Value = "Stakeholder S (1), Value v(1,i)"
$ws.Range("E1").Value = "Stakeholder S (1), Urgency u(1,i)"
$ws.Range("F1").Value = "Stakeholder S (2), Value v(2,i)"
$ws.Range("G1").Value = "Stakeholder S (2), Urgency u(2,i)"

# Match the author's saved selection.
$null = $ws.Range("B6").Select()
